# Weekly update: a new Mango price record (week of 2022-08-21, serial 44783)
# was reported for the Macroferia Regional de Talca market. It is inserted
# as the new row 126, pushing the existing rows 126-140 down to 127-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 126 - everything at/after 126
# shifts down by one (old 126 -> 127, ..., old 140 -> 141).
$ws.Rows.Item(126).Insert()

# The descriptive / categorical columns are identical for every record of
# this product+market, so copy them from the row that was just pushed down
# (now row 127, i.e. the former row 126) into the newly blank row 126.
$ws.Range("A126").Value2 = $ws.Range("A127").Value2()
$ws.Range("B126").Value2 = $ws.Range("B127").Value2()
$ws.Range("C126").Value2 = $ws.Range("C127").Value2()
$ws.Range("E126").Value2 = $ws.Range("E127").Value2()
$ws.Range("F126").Value2 = $ws.Range("F127").Value2()
$ws.Range("G126").Value2 = $ws.Range("G127").Value2()
$ws.Range("H126").Value2 = $ws.Range("H127").Value2()
$ws.Range("I126").Value2 = $ws.Range("I127").Value2()
$ws.Range("J126").Value2 = $ws.Range("J127").Value2()
$ws.Range("K126").Value2 = $ws.Range("K127").Value2()
$ws.Range("L126").Value2 = $ws.Range("L127").Value2()
$ws.Range("Q126").Value2 = $ws.Range("Q127").Value2()
$ws.Range("T126").Value2 = $ws.Range("T127").Value2()

# New values specific to this week's record.
$ws.Range("D126").Value2 = 44783
$ws.Range("M126").Value2 = 248
$ws.Range("N126").Value2 = 9000
$ws.Range("O126").Value2 = 9000
$ws.Range("P126").Value2 = 9000
$ws.Range("R126").Value2 = "México"
$ws.Range("S126").Value2 = 2250
